$d = $word.ActiveDocument

# 1. Fix the typo "BSIC" -> remove it so the heading reads "JUNIT TESTING-------"
$d.Content.Find.Execute("JUNIT BSIC TESTING-------", $true, $false, $false, $false, $false,
                         $true, 1, $false, "JUNIT TESTING-------", 2)

# 2. Mark the three picture-containing paragraphs as NoProofing so Word
#    writes <w:rPr><w:noProof/></w:rPr> on their runs (as happens when a
#    picture is re-pasted / regenerated by Word).
$pictureParaIndexes = @(2, 14, 19)
foreach ($idx in $pictureParaIndexes) {
    $para = $d.Paragraphs($idx)
    $para.Range.NoProofing = $true
}
